# Applies the diff:
#  - Adds w:proofErr spellStart/spellEnd and gramStart/gramEnd markers
#    (splitting some runs) around several tokens.
#  - Changes "randomForest_model.RDS" -> "predictive_model.RDS" (split into
#    two runs: "predictive" and "_model.RDS").
#  - Splits some "X -  a unique number key" style sentences with gramStart/End.
#  - Appends " with the predictions in the last column" after "predictions.csv".
#
# Because Range.InsertXML (when the fragment contains a <w:p> element)
# replaces the whole paragraph that the range lives in, each helper below
# rebuilds a complete paragraph (pPr + runs + proofErr) and reinserts it in
# place of the original paragraph text.

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml {
    param($Paragraph, $InnerXml)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wordNs + '><w:body>' +
        $InnerXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $result = $Paragraph.Range.InsertXML($pkg)
}

# --- Paragraph 2: "Both R object models output from topic_model_training.R"
$p2 = $d.Paragraphs.Item(2)
$xml2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Both R object models output from </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>topic_model_</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>training.R</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p2 $xml2

# --- Paragraph 3: "topic_model.RDS"
$p3 = $d.Paragraphs.Item(3)
$xml3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr>' +
    '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>topic_model.RDS</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p3 $xml3

# --- Paragraph 4: "randomForest_model.RDS" -> "predictive_model.RDS"
$p4 = $d.Paragraphs.Item(4)
$xml4 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr>' +
    '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>predictive</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>_model.RDS</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p4 $xml4

# --- Paragraph 7: "X  -  a unique number key"
$p7 = $d.Paragraphs.Item(7)
$xml7 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="6"/></w:numPr>' +
    '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">X </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> -</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve">  a unique number key</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p7 $xml7

# --- Paragraph 13: "The topic model -  R Object output from topic_model_training.R"
$p13 = $d.Paragraphs.Item(13)
$xml13 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="7"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">The topic model </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>-  R</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Object output from </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>topic_model_training.R</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p13 $xml13

# --- Paragraph 14: "The predictive model -  R Object output from topic_model_training.R"
$p14 = $d.Paragraphs.Item(14)
$xml14 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="7"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">The </w:t></w:r>' +
    '<w:r><w:t>predictive</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> model </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>-  R</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Object output from </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>topic_model_training.R</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p14 $xml14

# --- Paragraph 21: "The predictions are appended to the new report data set and
#     it is output as csv predictions.csv"
$p21 = $d.Paragraphs.Item(21)
$xml21 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">The predictions are appended to the new report data </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>set</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and it is output as csv </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>predictions.csv</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>with the predictions in the last column</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p21 $xml21

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
